$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'MSG: None

MSG: The decision about which movie to show on Friday resulted in no conclusion.
'
$ws.Range("C3").Value = 'MSG: None

MSG: The decision has been recorded to show "Barbie" on Friday.
'
$ws.Range("C4").Value = 'MSG: None

MSG: The rights to both movies have been acquired.
'
$ws.Range("C5").Value = 'MSG: None

MSG: The decision has been recorded, and the movie "Barbie" will be acquired for showing on Friday.
'
$ws.Range("C6").Value = 'MSG: None

MSG: The decision cannot be made regarding which movie to show on Friday.
'
$ws.Range("C7").Value = 'MSG: None

MSG: The decision-making process ended without a clear decision on which movie to show on Friday.
'
$ws.Range("C8").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Range("C9").Value = 'MSG: None

MSG: The decision has concluded with no choice of a movie for Friday.
'
$ws.Range("C10").Value = 'MSG: None

MSG: The decision has been recorded, and the rights for "Barbie" have been successfully acquired for the showing on Friday.
'
$ws.Range("D10").Value = 'Barbie_was_selected, '
$ws.Range("C11").Value = 'MSG: None

MSG: The decision to acquire the rights for "Barbie" has been made.
'
$ws.Range("C12").Value = 'MSG: None

MSG: The decision has been made, and there will be no movie shown on Friday.
'
$ws.Range("C13").Value = 'MSG: None

MSG: The decision process has concluded without selecting a movie for Friday.
'
$ws.Range("C14").Value = 'MSG: None

MSG: I have recorded the decision as "no_decision," indicating that no agreement was reached about the movie for Friday.
'
$ws.Range("D14").Value = 'no_decision, '
$ws.Range("C15").Value = 'MSG: None

MSG: The decision regarding Friday''s movie cannot be made at this time.
'
$ws.Range("C16").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights to "Barbie."
'
$ws.Range("C17").Value = 'MSG: None

MSG: The decision process concluded without a choice for Friday’s movie, so the no-decision function has been executed.
'
$ws.Range("C18").Value = 'MSG: None

MSG: A decision regarding the movie to show on Friday could not be reached, so I have recorded that as no decision made.
'
$ws.Range("C19").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Range("C20").Value = 'MSG: None

MSG: The decision process has concluded without selecting a movie for Friday.
'
$ws.Range("C21").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie will be shown on Friday as the committee did not reach a consensus.
'
$ws.Range("C22").Value = 'MSG: None

MSG: The conversation has ended without making a decision about what movie to show on Friday, so I will call the no_decision function.
'
$ws.Range("C23").Value = 'MSG: None

MSG: The decision has been recorded as no decision on which movie to show on Friday.
'
$ws.Range("C24").Value = 'MSG: None

MSG: The decision about which movie to show on Friday has not been made.
'
$ws.Range("C25").Value = 'MSG: None

MSG: It appears that no decision about Friday''s movie was made. Therefore, I will call the `no_decision` function to reflect this outcome.
'
$ws.Range("D25").Value = 'no_decision, '
$ws.Range("C26").Value = 'MSG: None

MSG: No decision about Friday''s movie was made in this meeting.
'
$ws.Range("C27").Value = 'MSG: None

MSG: No decision was reached regarding the movie to be shown on Friday.
'
$ws.Range("C28").Value = 'MSG: None

MSG: The rights to both movies have been acquired successfully.
'
$ws.Range("C29").Value = 'MSG: None

MSG: The decision to acquire the rights for the movie "Barbie" has been made.
'
$ws.Range("D29").Value = 'Barbie_was_selected, '
$ws.Range("C30").Value = 'MSG: None

MSG: The rights to both movies have been acquired.
'
$ws.Range("D30").Value = 'both_movies, '
$ws.Range("C31").Value = 'MSG: None

MSG: Based on the information provided, it appears that the committee was unable to reach a decision about what movie to show on Friday. Therefore, the appropriate action is to call the no_decision function.
'
$ws.Range("C32").Value = 'MSG: None

MSG: The decision has been recorded as "no decision" regarding the movie for Friday.
'
$ws.Range("C33").Value = 'MSG: None

MSG: The decision about which movie to show on Friday has not been made.
'
$ws.Range("C34").Value = 'MSG: None

MSG: The function has been executed, indicating that no decision was reached regarding the movie selection.
'
$ws.Range("C35").Value = 'MSG: None

MSG: The conversation ended without a decision on what movie to show on Friday. Therefore, the appropriate action was to call the no_decision function.
'
$ws.Range("C36").Value = 'MSG: None

MSG: The decision-making process did not result in a selection for Friday''s movie.
'
$ws.Range("C37").Value = 'MSG: None

MSG: The committee did not arrive at a decision regarding which movie to show on Friday.
'
$ws.Range("C38").Value = 'MSG: None

MSG: The decision has been recorded as no movie was selected.
'
$ws.Range("D38").Value = 'no_decision, '
$ws.Range("C39").Value = 'MSG: None

MSG: I have recorded the decision as no decision was made regarding the movie to be shown on Friday.
'
$ws.Range("C40").Value = 'MSG: None

MSG: I apologize for the confusion earlier. Since it seems there was no explicit decision reached regarding acquiring the rights for both movies or confirming one movie to show, I will proceed to call the function indicating that no decision was made about Friday’s movie.
```no_decision```
'
$ws.Range("D40").Value = 'both_movies, '
$ws.Range("C41").Value = 'MSG: None

MSG: I have recorded the decision that no movie will be shown on Friday.
'
$ws.Range("C42").Value = 'MSG: None

MSG: The decision to acquire the rights to "Barbie" has been successfully recorded.
'
$ws.Range("C43").Value = 'MSG: None

MSG: The decision has been recorded with no selection made for Friday''s movie.
'
$ws.Range("C44").Value = 'MSG: None

MSG: The decision to acquire the rights to show "Barbie" has been successfully recorded.
'
$ws.Range("C45").Value = 'MSG: None

MSG: The decision to acquire the rights for the movie "Barbie" has been recorded.
'
$ws.Range("C47").Value = 'MSG: None

MSG: No movie was chosen for Friday.
'
$ws.Range("C48").Value = 'MSG: None

MSG: It appears that the rights to both movies have been acquired.
'
$ws.Range("D48").Value = 'both_movies, '
$ws.Range("C49").Value = 'MSG: None

MSG: The decision regarding the movie to show on Friday resulted in no agreement being reached.
'
$ws.Range("C50").Value = 'MSG: None

MSG: The decision has been recorded, and the rights to "Oppenheimer" will be acquired.
'
$ws.Range("C51").Value = 'MSG: None

MSG: The decision about which movie to show on Friday remains unresolved.
'
$ws.Range("C52").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie choice was made for Friday.
'
$ws.Range("C53").Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday.
'
$ws.Range("C54").Value = 'MSG: None

MSG: The decision-making process did not result in an agreement on which movie to show, so no decision has been made regarding Friday''s movie.
'
$ws.Range("C55").Value = 'MSG: None

MSG: The rights to "Barbie" have been successfully acquired for Friday''s showing.
'
$ws.Range("C56").Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie."
'
$ws.Range("C57").Value = 'MSG: None

MSG: The decision to acquire the rights for both movies has been successfully recorded.
'
$ws.Range("C58").Value = 'MSG: None

MSG: The decision has been recorded as "no decision" regarding which movie to show on Friday, and no further action will be taken.
'
$ws.Range("C59").Value = 'MSG: None

MSG: The committee did not arrive at a decision about the movie to be shown on Friday. Therefore, no acquisition of movie rights is necessary.
'
$ws.Range("C60").Value = 'MSG: None

MSG: The decision to select a movie for Friday has resulted in no definitive choice being made.
'
$ws.Range("C61").Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not made.
'
$ws.Range("C62").Value = 'MSG: None

MSG: The decision has been recorded as no movie was selected in the meeting.
'
$ws.Range("D62").Value = 'no_decision, '
$ws.Range("C63").Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie" for screening on Friday.
'
$ws.Range("C64").Value = 'MSG: None

MSG: The decision has been recorded as "no decision."
'
$ws.Range("C66").Value = 'MSG: None

MSG: The decision was made that there would be no selected movie for Friday.
'
$ws.Range("C67").Value = 'MSG: None

MSG: The decision to show "Barbie" has been recorded successfully.
'
$ws.Range("C68").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie was selected for Friday.
'
